{"js": "// Update mvabund 1-year Q2 ANOVA table p-values and F-values after AIC\n// comparison / finalizing 1-year models.\nconst replacements = [\n  [\"0.997\", \"0.996\"],\n  [\"0.938\", \"0.939\"],\n  [\"0.577\", \"0.741\"],\n  [\"0.540\", \"0.597\"],\n  [\"5.651\", \"7.134\"],\n  [\"0.787\", \"0.781\"],\n];\n\nfor (const [oldVal, newVal] of replacements) {\n  const results = context.document.body.search(oldVal, { matchCase: true, matchWholeWord: false });\n  results.load(\"items/text\");\n  await context.sync();\n\n  if (results.items.length !== 1) {\n    throw new Error(`Expected exactly 1 match for \"${oldVal}\", found ${results.items.length}`);\n  }\n\n  results.items[0].insertText(newVal, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Update mvabund 1-year Q2 ANOVA table p-values and F-values after AIC\n# comparison / finalizing 1-year models.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{Old = \"0.997\"; New = \"0.996\"},\n    @{Old = \"0.938\"; New = \"0.939\"},\n    @{Old = \"0.577\"; New = \"0.741\"},\n    @{Old = \"0.540\"; New = \"0.597\"},\n    @{Old = \"5.651\"; New = \"7.134\"},\n    @{Old = \"0.787\"; New = \"0.781\"}\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $r.Old\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $r.New\n    $find.Execute(\n        $r.Old,   # FindText\n        $false,   # MatchCase\n        $true,    # MatchWholeWord\n        $false,   # MatchWildcards\n        $false,   # MatchSoundsLike\n        $false,   # MatchAllWordForms\n        $true,    # Forward\n        1,        # Wrap (wdFindContinue)\n        $false,   # Format\n        $r.New,   # ReplaceWith\n        2         # Replace (wdReplaceAll)\n    )\n}\n"}
